$wb = $excel.ActiveWorkbook

# --- Petition sheet (Petition / sheet4.xml) -------------------------------
$ws = $wb.Worksheets.Item("Petition")

# M5 changes from "Paper filing" to "n/a"
$ws.Range("M5").Value2 = "n/a"

# Duplicate row 6 into a brand-new row 7 (keeps values/styles identical to row 6)
$ws.Range("A6:V6").Copy($ws.Range("A7:V7"))

# Row 7 differs from row 6 in a few cells:
$ws.Range("C7").Value2 = 5
$ws.Range("G7").Value2 = "n/a"
$ws.Range("M7").Value2 = "Paper filing"

# --- CourtGenerateDocument sheet (sheet6.xml) -----------------------------
# Selection moves from M13 to D5 (this also keeps that sheet not the active one,
# as long as we re-activate Petition afterwards)
$ws2 = $wb.Worksheets.Item("CourtGenerateDocument")
$ws2.Range("D5").Select()

# --- Make Petition the active/selected sheet, with the selection on V7 ---
$ws.Activate()
$ws.Range("V7").Select()
